$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "b.md" has been handed off again; this report now reflects that new
# handoff: Status flips from "Handed back: in sync with en-US" to
# "Ready for handoff", a fresh handoff xliff + timestamp are recorded, and
# an error detail note about the stale handback is attached.
# ---------------------------------------------------------------------------

# --- Overview sheet: row for b.md (row 3) ---------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-22 22:37:12"

# --- zh-cn sheet: row for b.md (row 3) -------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# "Content Duplicate" flips to text "False"; copy from A2's sibling cell (F2,
# already literal text "False") so it lands as text, not a boolean.
$wsZhCn.Range("F2").Copy()
$wsZhCn.Range("F3").PasteSpecial(-4163)
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-22 22:37:07"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1f6c035d3366bdc4019a05e8855a4d560b2ca59f/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea25a580b8f09710557a552d902ffd1b24282c13/e2e/b.md."
# widen the "Error Detail" column to fit the new message
$zhCnRefWidth = $wsZhCn.Columns.Item(7).ColumnWidth()
$wsZhCn.Columns.Item(16).ColumnWidth = $zhCnRefWidth

# --- de-de sheet: row for b.md (row 3) -------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
# "Content Duplicate" flips to text "False"; copy from A2's sibling cell (F2,
# already literal text "False") so it lands as text, not a boolean.
$wsDeDe.Range("F2").Copy()
$wsDeDe.Range("F3").PasteSpecial(-4163)
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-22 22:37:12"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1f6c035d3366bdc4019a05e8855a4d560b2ca59f/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea25a580b8f09710557a552d902ffd1b24282c13/e2e/b.md."
# widen the "Error Detail" column to fit the new message
$deDeRefWidth = $wsDeDe.Columns.Item(7).ColumnWidth()
$wsDeDe.Columns.Item(16).ColumnWidth = $deDeRefWidth
